$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new data row (row 89) to the Facebook activity log.
$ws.Range("A89").Value = 44188                              # A: Date
$ws.Range("B89").Value = 0.46111111111111108                # B: Time
$ws.Range("C89").Value = "Friends of Friends"                # C: Privacy
$ws.Range("D89").Value = "Secretary status unlocked!"        # D: Post
$ws.Range("E89").Value = "10107916485805899"                 # E: Text
$ws.Range("F89").Value = 0                                   # F: Like
$ws.Range("G89").Value = 0                                   # G: Love
$ws.Range("H89").Value = 0                                   # H: Wow
$ws.Range("I89").Value = 1                                   # I: Haha
$ws.Range("J89").Value = 0                                   # J: Sad
$ws.Range("K89").Value = 0                                   # K: Care
$ws.Range("L89").Value = 0                                   # L: Angry
$ws.Range("M89").Value = 1                                   # M: Comments
$ws.Range("N89").Value = "10107916485805899"                 # N: Reference
$ws.Range("O89").Value = "Rohan Lewis"                        # O: Secretary

# Reflect the cursor/selection move that accompanied the new row.
$ws.Range("N93").Select() | Out-Null
